# Using CaseUtil and RestUtil: the ApiInfoSheet header row gains inline
# (parenthetical) descriptions matching the style already used on the
# Case2 sheet header row.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ApiInfoSheet")
$ws2 = $wb.Worksheets.Item("Case2")

# Reword the ApiInfoSheet header row to match CaseUtil/RestUtil's expected
# column labels (with Chinese field descriptions, like Case2's header).
$ws1.Range("A1").Value = "ApiId(接口编号)"
$ws1.Range("B1").Value = "ApiName(接口名称)"
$ws1.Range("C1").Value = "Type(接口提交方式)"
$ws1.Range("D1").Value = "Url(接口地址)"

# Widen column C on ApiInfoSheet to fit the new "Type(接口提交方式)" header.
$ws1.Columns.Item(3).ColumnWidth = 22.166666666666668

# Update the remembered selections on both sheets and flip which tab is
# active (ApiInfoSheet becomes the active/front tab, cursor on C2; Case2's
# remembered selection moves to D13).
$null = $ws2.Range("D13").Select()
$null = $ws1.Activate()
$null = $ws1.Range("C2").Select()
